# Auto-generated edit script applying the Ultima_Profits market-data refresh
# (prices/profits columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1566.5333
$ws.Range("I6").Value = 1208.1666
$ws.Range("J6").Value = 3000
$ws.Range("K6").Value = 3624.4998
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = -3512.4998
$ws.Range("N6").Value = -9224

$ws.Range("H9").Value = 115.8
$ws.Range("I9").Value = 150
$ws.Range("J9").Value = 93
$ws.Range("K9").Value = 150
$ws.Range("L9").Value = 93
$ws.Range("M9").Value = 19
$ws.Range("N9").Value = -431

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""

$ws.Range("H29").Value = 140
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""

$ws.Range("H38").Value = 1429468.4
$ws.Range("I38").Value = 96.666664
$ws.Range("J38").Value = 2501497.2
$ws.Range("K38").Value = 289.999992
$ws.Range("L38").Value = 7504491.600000001
$ws.Range("M38").Value = 82.00000799999998
$ws.Range("N38").Value = -7505235.600000001

$ws.Range("H40").Value = 1442.8572
$ws.Range("J40").Value = 1633.3334
$ws.Range("L40").Value = 1633.3334
$ws.Range("N40").Value = -1983.3334

$ws.Range("H58").Value = 445.9798
$ws.Range("I58").Value = 220.25
$ws.Range("J58").Value = 477.11493
$ws.Range("K58").Value = 660.75
$ws.Range("L58").Value = 1431.34479
$ws.Range("M58").Value = -510.75
$ws.Range("N58").Value = -1731.34479

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""

$ws.Range("H137").Value = 1015.5833
$ws.Range("I137").Value = 1097
$ws.Range("J137").Value = 934.1667
$ws.Range("K137").Value = 3291
$ws.Range("L137").Value = 2802.5001
$ws.Range("M137").Value = -741
$ws.Range("N137").Value = -7902.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3101.23
$ws.Range("I32").Value = 3073.9697
$ws.Range("K32").Value = 3073.9697
$ws.Range("M32").Value = -2786.9697

$ws.Range("H63").Value = 166670130
$ws.Range("I63").Value = 166670130
$ws.Range("K63").Value = 166670130
$ws.Range("M63").Value = -166669444

$ws.Range("H66").Value = 166670130
$ws.Range("I66").Value = 166670130
$ws.Range("K66").Value = 833350650
$ws.Range("M66").Value = -833347218

$ws.Range("H122").Value = 4783.385
$ws.Range("I122").Value = 5417.25
$ws.Range("K122").Value = 16251.75
$ws.Range("M122").Value = -13801.75

$ws.Range("H132").Value = 7355191
$ws.Range("I132").Value = 10418440
$ws.Range("K132").Value = 31255320
$ws.Range("M132").Value = -31252790

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1900.25
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""

$ws.Range("H35").Value = 39796.668
$ws.Range("J35").Value = 39796.668
$ws.Range("L35").Value = 39796.668
$ws.Range("N35").Value = -40416.668

$ws.Range("H134").Value = 3862.566
$ws.Range("I134").Value = 2706.6155
$ws.Range("J134").Value = 4975.7036
$ws.Range("K134").Value = 8119.8465
$ws.Range("L134").Value = 14927.1108
$ws.Range("M134").Value = -5584.8465
$ws.Range("N134").Value = -19997.1108

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2141.0625
$ws.Range("I58").Value = 1318.421
$ws.Range("J58").Value = 3343.3845
$ws.Range("K58").Value = 1318.421
$ws.Range("L58").Value = 3343.3845
$ws.Range("M58").Value = -1115.421
$ws.Range("N58").Value = -3749.3845

$ws.Range("H136").Value = 2141.0625
$ws.Range("I136").Value = 1318.421
$ws.Range("J136").Value = 3343.3845
$ws.Range("K136").Value = 3955.263
$ws.Range("L136").Value = 10030.1535
$ws.Range("M136").Value = -1405.263
$ws.Range("N136").Value = -15130.1535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1200.7755
$ws.Range("I131").Value = 651.38464
$ws.Range("J131").Value = 1399.1666
$ws.Range("K131").Value = 1954.15392
$ws.Range("L131").Value = 4197.4998
$ws.Range("M131").Value = 3085.84608
$ws.Range("N131").Value = -14277.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 136.375
$ws.Range("I2").Value = 88.2
$ws.Range("K2").Value = 88.2
$ws.Range("M2").Value = 24.8

$ws.Range("H43").Value = 6000

$ws.Range("H57").Value = 11485.571
$ws.Range("J57").Value = 29000
$ws.Range("L57").Value = 29000
$ws.Range("N57").Value = -30640

$ws.Range("H80").Value = 15875717
$ws.Range("I80").Value = 30305396
$ws.Range("J80").Value = 3070
$ws.Range("K80").Value = 30305396
$ws.Range("L80").Value = 3070
$ws.Range("M80").Value = -30304398
$ws.Range("N80").Value = -5066

$ws.Range("H83").Value = 15875717
$ws.Range("I83").Value = 30305396
$ws.Range("J83").Value = 3070
$ws.Range("K83").Value = 151526980
$ws.Range("L83").Value = 15350
$ws.Range("M83").Value = -151521988
$ws.Range("N83").Value = -25334

$ws.Range("H132").Value = 5220.811
$ws.Range("I132").Value = 5898.5
$ws.Range("J132").Value = 3112.4443
$ws.Range("K132").Value = 17695.5
$ws.Range("L132").Value = 9337.332900000001
$ws.Range("M132").Value = -15165.5
$ws.Range("N132").Value = -14397.3329

$ws.Range("H136").Value = 23764.637
$ws.Range("J136").Value = 21208.5
$ws.Range("L136").Value = 63625.5
$ws.Range("N136").Value = -68725.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 21049
$ws.Range("I54").Value = 10345
$ws.Range("J54").Value = 23725
$ws.Range("K54").Value = 10345
$ws.Range("L54").Value = 23725
$ws.Range("M54").Value = -9825
$ws.Range("N54").Value = -24765

$ws.Range("H62").Value = 5084
$ws.Range("I62").Value = 4730.769
$ws.Range("J62").Value = 5594.222
$ws.Range("K62").Value = 4730.769
$ws.Range("L62").Value = 5594.222
$ws.Range("M62").Value = -4106.769
$ws.Range("N62").Value = -6842.222

$ws.Range("H65").Value = 5084
$ws.Range("I65").Value = 4730.769
$ws.Range("J65").Value = 5594.222
$ws.Range("K65").Value = 23653.845
$ws.Range("L65").Value = 27971.11
$ws.Range("M65").Value = -20533.845
$ws.Range("N65").Value = -34211.11

$ws.Range("H75").Value = 38700
$ws.Range("J75").Value = 38700
$ws.Range("L75").Value = 38700
$ws.Range("N75").Value = -40572

$ws.Range("H78").Value = 38700
$ws.Range("J78").Value = 38700
$ws.Range("L78").Value = 116100
$ws.Range("N78").Value = -125460

$ws.Range("H81").Value = 693.1875
$ws.Range("J81").Value = 2133
$ws.Range("L81").Value = 4266
$ws.Range("N81").Value = -6388

$ws.Range("H84").Value = 693.1875
$ws.Range("J84").Value = 2133
$ws.Range("L84").Value = 21330
$ws.Range("N84").Value = -31938
